$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp label in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 18 de Abril de 2020 a las 15:22"

# Country data refresh + 3 re-sorted country rows (Serbia/Singapur swap,
# Ghana re-inserted before Republica de Chipre shifting rows down, Cabo Verde/Zambia swap)
$ws.Cells.Item(8, 1).Value = "Alemania"
$ws.Cells.Item(8, 2).Value = 141968
$ws.Cells.Item(8, 3).Value = 571
$ws.Cells.Item(8, 4).Value = 85400
$ws.Cells.Item(8, 5).Value = 52191
$ws.Cells.Item(8, 6).Value = 5013
$ws.Cells.Item(8, 7).Value = 25
$ws.Cells.Item(8, 8).Value = 4377

$ws.Cells.Item(9, 1).Value = "Reino Unido"
$ws.Cells.Item(9, 2).Value = 114217
$ws.Cells.Item(9, 3).Value = 5525
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 98409
$ws.Cells.Item(9, 6).Value = 1559
$ws.Cells.Item(9, 7).Value = 888
$ws.Cells.Item(9, 8).Value = 15464

$ws.Cells.Item(42, 1).Value = "Serbia"
$ws.Cells.Item(42, 2).Value = 5994
$ws.Cells.Item(42, 3).Value = 304
$ws.Cells.Item(42, 4).Value = 637
$ws.Cells.Item(42, 5).Value = 5240
$ws.Cells.Item(42, 6).Value = 120
$ws.Cells.Item(42, 7).Value = 7
$ws.Cells.Item(42, 8).Value = 117

$ws.Cells.Item(43, 1).Value = "Singapur"
$ws.Cells.Item(43, 2).Value = 5992
$ws.Cells.Item(43, 3).Value = 942
$ws.Cells.Item(43, 4).Value = 708
$ws.Cells.Item(43, 5).Value = 5273
$ws.Cells.Item(43, 6).Value = 22
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 11

$ws.Cells.Item(70, 1).Value = "Uzbekistan"
$ws.Cells.Item(70, 2).Value = 1450
$ws.Cells.Item(70, 3).Value = 45
$ws.Cells.Item(70, 4).Value = 194
$ws.Cells.Item(70, 5).Value = 1252
$ws.Cells.Item(70, 6).Value = 8
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 4

$ws.Cells.Item(72, 1).Value = "Azerbaiyan"
$ws.Cells.Item(72, 2).Value = 1373
$ws.Cells.Item(72, 3).Value = 33
$ws.Cells.Item(72, 4).Value = 590
$ws.Cells.Item(72, 5).Value = 765
$ws.Cells.Item(72, 6).Value = 21
$ws.Cells.Item(72, 7).Value = 3
$ws.Cells.Item(72, 8).Value = 18

$ws.Cells.Item(78, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(78, 2).Value = 1170
$ws.Cells.Item(78, 3).Value = 53
$ws.Cells.Item(78, 4).Value = 164
$ws.Cells.Item(78, 5).Value = 957
$ws.Cells.Item(78, 6).Value = 15
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 49

$ws.Cells.Item(86, 1).Value = "Ghana"
$ws.Cells.Item(86, 2).Value = 834
$ws.Cells.Item(86, 3).Value = 193
$ws.Cells.Item(86, 4).Value = 99
$ws.Cells.Item(86, 5).Value = 726
$ws.Cells.Item(86, 6).Value = 4
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 9

$ws.Cells.Item(87, 1).Value = "Republica de Chipre"
$ws.Cells.Item(87, 2).Value = 750
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = 77
$ws.Cells.Item(87, 5).Value = 661
$ws.Cells.Item(87, 6).Value = 8
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 12

$ws.Cells.Item(88, 1).Value = "Costa de Marfil"
$ws.Cells.Item(88, 2).Value = 742
$ws.Cells.Item(88, 3).Value = 54
$ws.Cells.Item(88, 4).Value = 220
$ws.Cells.Item(88, 5).Value = 516
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 6

$ws.Cells.Item(89, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(89, 2).Value = 732
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 76
$ws.Cells.Item(89, 5).Value = 654
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 2

$ws.Cells.Item(90, 1).Value = "Letonia"
$ws.Cells.Item(90, 2).Value = 712
$ws.Cells.Item(90, 3).Value = 30
$ws.Cells.Item(90, 4).Value = 88
$ws.Cells.Item(90, 5).Value = 619
$ws.Cells.Item(90, 6).Value = 5
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 5

$ws.Cells.Item(91, 1).Value = "Crucero"
$ws.Cells.Item(91, 2).Value = 712
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 4).Value = 644
$ws.Cells.Item(91, 5).Value = 55
$ws.Cells.Item(91, 6).Value = 7
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 13

$ws.Cells.Item(92, 1).Value = "Principado de Andorra"
$ws.Cells.Item(92, 2).Value = 696
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 191
$ws.Cells.Item(92, 5).Value = 470
$ws.Cells.Item(92, 6).Value = 17
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 35

$ws.Cells.Item(93, 1).Value = "Libano"
$ws.Cells.Item(93, 2).Value = 672
$ws.Cells.Item(93, 3).Value = 4
$ws.Cells.Item(93, 4).Value = 94
$ws.Cells.Item(93, 5).Value = 557
$ws.Cells.Item(93, 6).Value = 30
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 21

$ws.Cells.Item(94, 1).Value = "Costa Rica"
$ws.Cells.Item(94, 2).Value = 649
$ws.Cells.Item(94, 3).Value = 0
$ws.Cells.Item(94, 4).Value = 88
$ws.Cells.Item(94, 5).Value = 557
$ws.Cells.Item(94, 6).Value = 10
$ws.Cells.Item(94, 7).Value = 0
$ws.Cells.Item(94, 8).Value = 4

$ws.Cells.Item(124, 1).Value = "El Salvador"
$ws.Cells.Item(124, 2).Value = 190
$ws.Cells.Item(124, 3).Value = 13
$ws.Cells.Item(124, 4).Value = 43
$ws.Cells.Item(124, 5).Value = 140
$ws.Cells.Item(124, 6).Value = 1
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 7

$ws.Cells.Item(154, 1).Value = "Cabo Verde"
$ws.Cells.Item(154, 2).Value = 58
$ws.Cells.Item(154, 3).Value = 2
$ws.Cells.Item(154, 4).Value = 1
$ws.Cells.Item(154, 5).Value = 56
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 1

$ws.Cells.Item(155, 1).Value = "Zambia"
$ws.Cells.Item(155, 2).Value = 57
$ws.Cells.Item(155, 3).Value = 5
$ws.Cells.Item(155, 4).Value = 33
$ws.Cells.Item(155, 5).Value = 22
$ws.Cells.Item(155, 6).Value = 1
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 2

